$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.456.51"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.851.41"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.47"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6293"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07660"
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.82"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "2.508.68"
$ws.Range("E11").Value = "  +34.11%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07763"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.039"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6820"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001056"
$ws.Range("E15").Value = "  -4.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.54"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.198"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.459.71"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.50"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.512"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.62"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1386"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.438"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.377"
$ws.Range("E28").Value = "  +5.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.465"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05613"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.138"
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.062"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.847"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.168"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7007"
$ws.Range("E35").Value = "  -1.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.600"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01805"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.226.77"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.742"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.450"
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9077"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.13"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.28"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.219"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000119"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4026"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.071"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1156"
$ws.Range("E49").Value = "  +3.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.679"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05713"
$ws.Range("E51").Value = "  +0.02%  "
